$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B2")
$text = $cell.Value2
$newtext = $text.Replace("RME/", "")
$cell.Value = $newtext
$cell.WrapText = $true
$ws.Rows.Item(2).RowHeight = 409.6
$null = $ws.Range("E4").Select()
